$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.616.58"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "3.447.36"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'578.12"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'148.64"
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("E9").Value = "  +5.54%  "

$ws.Range("E10").Value = "  -0.88%  "

$ws.Range("D11").Value = "'0.414"
$ws.Range("E11").Value = "  +3.38%  "

$ws.Range("D12").Value = "4.039.42"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("E13").Value = "  -0.06%  "

$ws.Range("D14").Value = "'28.31"
$ws.Range("E14").Value = "  -5.19%  "

$ws.Range("D15").Value = "3.447.48"
$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("E16").Value = "  +0.77%  "

$ws.Range("D17").Value = "62.691.34"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "'6.40"
$ws.Range("E18").Value = "  +0.72%  "

$ws.Range("D19").Value = "'14.60"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("E20").Value = "  -2.76%  "

$ws.Range("D21").Value = "'386.82"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("D22").Value = "'0.568"
$ws.Range("E22").Value = "  +0.65%  "

$ws.Range("D23").Value = "'75.23"
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "3.584.92"
$ws.Range("E25").Value = "  -0.67%  "

$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("E27").Value = "  +2.68%  "

$ws.Range("D28").Value = "'7.66"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("E31").Value = "  -1.32%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  -3.66%  "

$ws.Range("D34").Value = "'23.22"
$ws.Range("E34").Value = "  -2.20%  "

$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  +1.84%  "

$ws.Range("D36").Value = "'1.63"
$ws.Range("E36").Value = "  +3.70%  "

$ws.Range("D37").Value = "'32.03"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("D39").Value = "'169.19"
$ws.Range("E39").Value = "  -0.32%  "

$ws.Range("D40").Value = "3.482.68"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("E41").Value = "  +0.77%  "

$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("D43").Value = "'42.70"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("D44").Value = "'4.38"
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("E45").Value = "  -1.68%  "

$ws.Range("E46").Value = "  -1.57%  "

$ws.Range("D47").Value = "2.567.23"
$ws.Range("E47").Value = "  -1.12%  "

$ws.Range("E48").Value = "  +1.82%  "

$ws.Range("E49").Value = "  +1.52%  "

$ws.Range("D50").Value = "'22.54"
$ws.Range("E50").Value = "  -3.71%  "

$ws.Range("E51").Value = "  +0.01%  "
